$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last refreshed" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 09:51"

# --- Country data refresh ---
# The underlying source data was re-pulled; several countries' totals grew
# enough to overtake their neighbour in the (descending, by "Casos totales")
# sorted list, so those two rows swap position with fresh numbers, and every
# row in between shifts down by one (keeping its old numbers). All other
# rows are simple in-place value updates.

# India (row 5)
$ws.Range("B5").Value = 4208645
$ws.Range("C5").Value = 6083
$ws.Range("E5").Value = 886505
$ws.Range("G5").Value = 24
$ws.Range("H5").Value = 71711

# Rusia (row 7)
$ws.Range("B7").Value = 1030690
$ws.Range("C7").Value = 5185
$ws.Range("D7").Value = 843277
$ws.Range("E7").Value = 169542
$ws.Range("G7").Value = 51
$ws.Range("H7").Value = 17871

# Singapur (row 53)
$ws.Range("B53").Value = 57044
$ws.Range("C53").Value = 22
$ws.Range("E53").Value = 684

# Barein (row 54)
$ws.Range("E54").Value = 4269
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 200

# Armenia (row 60)
$ws.Range("B60").Value = 44845
$ws.Range("C60").Value = 62
$ws.Range("D60").Value = 40121
$ws.Range("E60").Value = 3824
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 900

# Row 98: Tayikistan -> Hungria (Hungria overtakes Tayikistan, new data)
$ws.Range("A98").Value = "Hungria"
$ws.Range("B98").Value = 8963
$ws.Range("C98").Value = 576
$ws.Range("D98").Value = 3961
$ws.Range("E98").Value = 4377
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 625

# Row 99: Namibia -> Tayikistan (shifted down, old Tayikistan numbers)
$ws.Range("A99").Value = "Tayikistan"
$ws.Range("B99").Value = 8792
$ws.Range("D99").Value = 7581
$ws.Range("E99").Value = 1141
$ws.Range("H99").Value = 70

# Row 100: Gabon -> Namibia (shifted down, old Namibia numbers)
$ws.Range("A100").Value = "Namibia"
$ws.Range("B100").Value = 8685
$ws.Range("D100").Value = 3786
$ws.Range("E100").Value = 4810
$ws.Range("H100").Value = 89

# Row 101: Maldivas -> Gabon (shifted down, old Gabon numbers)
$ws.Range("A101").Value = "Gabon"
$ws.Range("B101").Value = 8601
$ws.Range("D101").Value = 7424
$ws.Range("E101").Value = 1124
$ws.Range("H101").Value = 53

# Row 102: Hungria -> Maldivas (shifted down, old Maldivas numbers)
$ws.Range("A102").Value = "Maldivas"
$ws.Range("B102").Value = 8584
$ws.Range("D102").Value = 5936
$ws.Range("E102").Value = 2619
$ws.Range("H102").Value = 29

# Row 106: Luxemburgo -> Zimbabue (Zimbabue overtakes Luxemburgo, new data)
$ws.Range("A106").Value = "Zimbabue"
$ws.Range("B106").Value = 7116
$ws.Range("C106").Value = 279
$ws.Range("D106").Value = 5373
$ws.Range("E106").Value = 1535
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 208

# Row 107: Zimbabue -> Luxemburgo (shifted down, old Luxemburgo numbers)
$ws.Range("A107").Value = "Luxemburgo"
$ws.Range("B107").Value = 6950
$ws.Range("D107").Value = 6126
$ws.Range("E107").Value = 700
$ws.Range("H107").Value = 124

# Letonia (row 157)
$ws.Range("B157").Value = 1429
$ws.Range("C157").Value = 1
$ws.Range("E157").Value = 207
